$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated ligand average/total expression values (same new TPM-derived value across rows 2-4)
$ws.Range("G2:G4").Value = 0.004257666666666667
$ws.Range("H2:H4").Value = 0.012773

# Row 2 (Osm / Il6st / ECs -> ECs) receptor + edge values
$ws.Range("M2").Value = 19.827687
$ws.Range("N2").Value = 59.483061
$ws.Range("O2").Value = 0.1538389073329896
$ws.Range("P2").Value = 0.1538389073329896
$ws.Range("Q2").Value = 0.08441968201700001
$ws.Range("R2").Value = 0.759777138153
$ws.Range("S2").Value = 0.1538389073329896
$ws.Range("T2").Value = 0.1538389073329896

# Row 3 (Osm / Il6st / ECs -> FAPs) edge + specificity values (M3/N3 unchanged)
$ws.Range("O3").Value = 0.6604253914664442
$ws.Range("P3").Value = 0.6604253914664441
$ws.Range("Q3").Value = 0.3624109304343334
$ws.Range("R3").Value = 3.261698373909
$ws.Range("S3").Value = 0.6604253914664442
$ws.Range("T3").Value = 0.6604253914664441

# Row 4 (Osm / Il6st / ECs -> MuSCs) receptor + edge values
$ws.Range("M4").Value = 23.93873833333333
$ws.Range("N4").Value = 71.816215
$ws.Range("O4").Value = 0.1857357012005663
$ws.Range("P4").Value = 0.1857357012005663
$ws.Range("Q4").Value = 0.1019231682438889
$ws.Range("R4").Value = 0.9173085141949999
$ws.Range("S4").Value = 0.1857357012005663
$ws.Range("T4").Value = 0.1857357012005663
